# Correction in SA algorithm and 746 logs
# Update Fitness (column C) values for run_27 log to reflect corrected
# fitness computation. Rows 2-17 -> 7721, rows 18-19 -> 7318,
# rows 20-130 -> 7293. Rows 131-252 already hold the corrected value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C17").Value = 7721
$ws.Range("C18:C19").Value = 7318
$ws.Range("C20:C130").Value = 7293
